# Generate Report for Handoff
# Updates the localization-status report: marks the pending translations as
# "Ready for handoff" and refreshes the handoff timestamps.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: per-language status columns + overall latest handoff date
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-24-11 14:24:28"

# zh-cn detail sheet: status + handoff datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-11 14:24:25"

# de-de detail sheet: status + handoff datetime
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-11 14:24:28"
